# Update the embedded build timestamp throughout the workbook.
# "February 03 2026 17.29.55 EST" -> "February 03 2026 18.05.36 EST"

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$wb = $excel.ActiveWorkbook

# --- "About" sheet ---
$aboutWs = $wb.Worksheets.Item("About")

$aboutWs.Range("A2").Value = $aboutWs.Range("A2").Value().Replace($oldStamp, $newStamp)
$aboutWs.Range("A6").Value = $aboutWs.Range("A6").Value().Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet ---
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 13; $row++) {
    $cell = $dataWs.Cells.Item($row, 19)  # column S
    $cell.Value = $cell.Value().Replace($oldStamp, $newStamp)
}
